# Insert two new data rows (223 and 224) into the daily-price log sheet,
# pushing the existing rows 223-244 down to 225-246.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the existing row-223 block.
$ws.Rows("223:224").Insert()

# --- Row 223: Tomate / Larga vida / Primera, $/caja 10 kilos ---
$ws.Range("A223").Value = 11
$ws.Range("B223").Value = 'Vega Monumental Concepción'
$ws.Range("C223").Value = 'Bíobío'
$ws.Range("D223").Value = 44461
$ws.Range("E223").Value = 8
$ws.Range("F223").Value = 100112020
$ws.Range("G223").Value = 'Tomate'
$ws.Range("H223").Value = 'Larga vida'
$ws.Range("I223").Value = 'Primera'
$ws.Range("J223").Value = 2000
$ws.Range("K223").Value = 7000
$ws.Range("L223").Value = 7500
$ws.Range("M223").Value = 7250
$ws.Range("N223").Value = '$/caja 10 kilos'
$ws.Range("O223").Value = 'Región de Arica y Parinacota'
$ws.Range("P223").Value = 725
$ws.Range("Q223").Value = 10
$ws.Range("R223").Value = 'Hortaliza'

# --- Row 224: Tomate / Larga vida / Segunda, $/caja 10 kilos ---
$ws.Range("A224").Value = 11
$ws.Range("B224").Value = 'Vega Monumental Concepción'
$ws.Range("C224").Value = 'Bíobío'
$ws.Range("D224").Value = 44461
$ws.Range("E224").Value = 8
$ws.Range("F224").Value = 100112020
$ws.Range("G224").Value = 'Tomate'
$ws.Range("H224").Value = 'Larga vida'
$ws.Range("I224").Value = 'Segunda'
$ws.Range("J224").Value = 1000
$ws.Range("K224").Value = 6500
$ws.Range("L224").Value = 6500
$ws.Range("M224").Value = 6500
$ws.Range("N224").Value = '$/caja 10 kilos'
$ws.Range("O224").Value = 'Región de Arica y Parinacota'
$ws.Range("P224").Value = 650
$ws.Range("Q224").Value = 10
$ws.Range("R224").Value = 'Hortaliza'
